# Trading update: 2026-02-18 10:54:12
# Two new MarketMaking trades (Trade #5 @10:53:57 DOWN 0.28, Trade #6 @10:54:03 UP 0.51)
# get appended to the "All Trades" log, and the "MarketMaking" strategy sheet
# (which tracks only the currently-open trades for that strategy) is rolled
# forward: the prior open trade (#4) is replaced by the new open trade (#5)
# in row 2, and the newly-opened trade (#6) is appended as row 3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "All Trades" sheet
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #4 (row 5) is no longer the latest open MarketMaking trade, so its
# "still live on the strategy sheet" fields get cleared out, matching the
# pattern already used by the other non-latest OPEN rows (2-4): Exit Price
# becomes a numeric 0 instead of blank, and the capital/slippage/confidence/
# entry-reason/duration fields become blank.
$allTrades.Cells.Item(5, 7).Value = 0
$allTrades.Cells.Item(5, 11).Value = ""
$allTrades.Cells.Item(5, 12).Value = ""
$allTrades.Cells.Item(5, 13).Value = ""
$allTrades.Cells.Item(5, 14).Value = ""
$allTrades.Cells.Item(5, 15).Value = ""
$allTrades.Cells.Item(5, 17).Value = ""

# New row 6: Trade #5
$allTrades.Cells.Item(6, 1).Value = 5
$allTrades.Cells.Item(6, 2).NumberFormat = "@"
$allTrades.Cells.Item(6, 2).Value = "2026-02-18"
$allTrades.Cells.Item(6, 3).Value = "10:53:57"
$allTrades.Cells.Item(6, 4).Value = "MarketMaking"
$allTrades.Cells.Item(6, 5).Value = "DOWN"
$allTrades.Cells.Item(6, 6).Value = 0.28
$allTrades.Cells.Item(6, 7).Value = ""
$allTrades.Cells.Item(6, 8).Value = "OPEN"
$allTrades.Cells.Item(6, 9).Value = 0
$allTrades.Cells.Item(6, 10).Value = 0
$allTrades.Cells.Item(6, 11).Value = 100
$allTrades.Cells.Item(6, 12).Value = 0
$allTrades.Cells.Item(6, 13).Value = 0
$allTrades.Cells.Item(6, 14).Value = 0.6
$allTrades.Cells.Item(6, 15).Value = "Normal spread capture: 194 bps"
$allTrades.Cells.Item(6, 16).Value = ""
$allTrades.Cells.Item(6, 17).Value = 0

# New row 7: Trade #6
$allTrades.Cells.Item(7, 1).Value = 6
$allTrades.Cells.Item(7, 2).NumberFormat = "@"
$allTrades.Cells.Item(7, 2).Value = "2026-02-18"
$allTrades.Cells.Item(7, 3).Value = "10:54:03"
$allTrades.Cells.Item(7, 4).Value = "MarketMaking"
$allTrades.Cells.Item(7, 5).Value = "UP"
$allTrades.Cells.Item(7, 6).Value = 0.51
$allTrades.Cells.Item(7, 7).Value = ""
$allTrades.Cells.Item(7, 8).Value = "OPEN"
$allTrades.Cells.Item(7, 9).Value = 0
$allTrades.Cells.Item(7, 10).Value = 0
$allTrades.Cells.Item(7, 11).Value = 100
$allTrades.Cells.Item(7, 12).Value = 0
$allTrades.Cells.Item(7, 13).Value = 0
$allTrades.Cells.Item(7, 14).Value = 0.6
$allTrades.Cells.Item(7, 15).Value = "Normal spread capture: 194 bps"
$allTrades.Cells.Item(7, 16).Value = ""
$allTrades.Cells.Item(7, 17).Value = 0

# ---------------------------------------------------------------
# "MarketMaking" sheet
# ---------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Row 2 rolls forward from Trade #4 to Trade #5 (new open position).
$mm.Cells.Item(2, 1).Value = 5
$mm.Cells.Item(2, 3).Value = "10:53:57"
$mm.Cells.Item(2, 6).Value = 0.28
$mm.Cells.Item(2, 15).Value = "Normal spread capture: 194 bps"

# New row 3: Trade #6, also open on the MarketMaking strategy.
$mm.Cells.Item(3, 1).Value = 6
$mm.Cells.Item(3, 2).NumberFormat = "@"
$mm.Cells.Item(3, 2).Value = "2026-02-18"
$mm.Cells.Item(3, 3).Value = "10:54:03"
$mm.Cells.Item(3, 4).Value = "MarketMaking"
$mm.Cells.Item(3, 5).Value = "UP"
$mm.Cells.Item(3, 6).Value = 0.51
$mm.Cells.Item(3, 7).Value = ""
$mm.Cells.Item(3, 8).Value = "OPEN"
$mm.Cells.Item(3, 9).Value = 0
$mm.Cells.Item(3, 10).Value = 0
$mm.Cells.Item(3, 11).Value = 100
$mm.Cells.Item(3, 12).Value = 0
$mm.Cells.Item(3, 13).Value = 0
$mm.Cells.Item(3, 14).Value = 0.6
$mm.Cells.Item(3, 15).Value = "Normal spread capture: 194 bps"
$mm.Cells.Item(3, 16).Value = ""
$mm.Cells.Item(3, 17).Value = 0
